# BeverageOrders/Orders.xlsx - apply the commit's data edits.
#
# Row 2 ("Chai") order quantity changed from 2900 to -50, so the failure
# note text is updated to match.
# Row 4 changed from "Sasquatch Ale" (qty 4) to "Cote de Blaye" (qty 5).
# Finally, the grid selection on the Orders sheet moves to A4:B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

$ws.Range("B2").Value = -50
$ws.Range("D2").Value = "Quantity '-50' was unavailable"

$ws.Range("A4").Value = "Cote de Blaye"
$ws.Range("B4").Value = 5

$ws.Activate()
$ws.Range("A4:B4").Select() | Out-Null
